$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-15 Friday" "2024-03-16 Saturday"

Replace-Text "390÷3=" "976÷2="
Replace-Text "815÷7=" "308÷4="
Replace-Text "873÷7=" "535÷7="
Replace-Text "935÷5=" "358÷8="
Replace-Text "176÷8=" "604÷7="

Replace-Text "172÷2=" "605÷2="
Replace-Text "743÷7=" "547÷4="
Replace-Text "769÷7=" "129÷9="
Replace-Text "970÷8=" "448÷2="
Replace-Text "389÷3=" "238÷9="

Replace-Text "163÷6=" "251÷5="
Replace-Text "596÷5=" "584÷8="
Replace-Text "324÷3=" "839÷8="
Replace-Text "811÷5=" "530÷2="
Replace-Text "330÷2=" "491÷8="

Replace-Text "650÷3=" "448÷8="
Replace-Text "525÷4=" "597÷8="
Replace-Text "483÷8=" "485÷5="
Replace-Text "487÷7=" "333÷8="
Replace-Text "896÷8=" "286÷2="

Replace-Text "797÷5=" "150÷9="
Replace-Text "733÷4=" "579÷8="
Replace-Text "219÷3=" "908÷7="
Replace-Text "843÷8=" "106÷7="
Replace-Text "435÷7=" "253÷7="
